$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# text explicitly, otherwise Excel auto-converts the assigned string into
# a numeric value (e.g. "247.57" -> 247.57, dropping the original text
# representation). We flip the format to Text, assign the value, then
# restore the cell to the default (Normal) style so no extra formatting
# is left behind.

$ws.Range('D2').Value = '34.735.56'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').Value = '1.874.36'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  -0.89%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.689'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.41%  '
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.92'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.347'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0735'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0970'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').Value = '2.149.76'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.714'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').Value = '1.864.29'
$ws.Range('E17').Value = '  -3.01%  '
$ws.Range('D18').Value = '34.767.15'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('E25').Value = '  +4.05%  '
$ws.Range('E26').Value = '  -5.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.127'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.95%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0577'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.830'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.19%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.64'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -16.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '97.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.87'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0660'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0210'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  -4.73%  '
$ws.Range('D45').Value = '1.283.79'
$ws.Range('E45').Value = '  -4.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.47%  '
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0780'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.12%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.47%  '
